# Commit "commit the code upto 19/06/2025":
#  - add a new worksheet "adminusers" after the existing sheets
#    (becomes sheetId=3 / rId3 and the active tab)
#  - give it two small rows of shared-string data (A1 header +
#    A2/B2 repeated value), selection left on B2
#  - because the new sheet becomes the active one, the previously
#    active "loginpage" sheet is no longer the selected tab

$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the workbook (after the last
# existing worksheet) and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "adminusers"

# Populate the new sheet's data.
$ws.Range("A1").Value = "UsernamePassword"
$ws.Range("A2").Value = "pbv"
$ws.Range("B2").Value = "pbv"

# Leave the selection on B2, matching the saved view state.
$null = $ws.Range("B2").Select()
